# Update "想去人数" (want-to-go count) and "最低票价" (lowest price) figures
# for a handful of events that appear on both the "展览" (Exhibitions) sheet
# and the "全部类型" (All types) sheet, matching the refreshed data pull.

$wb = $excel.ActiveWorkbook

# Sheet 1 = "展览"
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("G3").Value  = 68
$wsExhibit.Range("F5").Value  = 460
$wsExhibit.Range("F6").Value  = 89
$wsExhibit.Range("F9").Value  = 6906
$wsExhibit.Range("F10").Value = 169
$wsExhibit.Range("F15").Value = 1120
$wsExhibit.Range("F16").Value = 16388
$wsExhibit.Range("F17").Value = 6
$wsExhibit.Range("F18").Value = 1611
$wsExhibit.Range("F20").Value = 342
$wsExhibit.Range("F23").Value = 11457
$wsExhibit.Range("F25").Value = 1097
$wsExhibit.Range("F26").Value = 4511

# Sheet 4 = "全部类型" (same underlying events, different row offsets)
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("G3").Value  = 68
$wsAll.Range("F5").Value  = 460
$wsAll.Range("F6").Value  = 89
$wsAll.Range("F10").Value = 6906
$wsAll.Range("F11").Value = 169
$wsAll.Range("F17").Value = 1120
$wsAll.Range("F18").Value = 16388
$wsAll.Range("F19").Value = 6
$wsAll.Range("F20").Value = 1611
$wsAll.Range("F22").Value = 342
$wsAll.Range("F27").Value = 11457
$wsAll.Range("F29").Value = 1097
$wsAll.Range("F30").Value = 4511
